$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 17 data
$ws.Range("A17").Value = 16
$ws.Range("B17").Value = 1.2568981481481483
$ws.Range("B17").NumberFormat = $ws.Range("B2").NumberFormat
$ws.Range("C17").Formula = "=SUM(B2:B17)+1.2708333333"
$ws.Range("C17").NumberFormat = $ws.Range("C2").NumberFormat
$ws.Range("D17").Value = "El Dragón: el regreso de un guerrero (Audiovisual, Spanish, New):36; [La vida de los africanos nómadas ganaderos | Mundari: Sudán del Sur 🇸🇸](https://youtu.be/QX4SIUwYOHE) (Audiovisual, Spanish, New):39; [10 COSAS que NO TIENEN SENTIDO en los ANIMES](https://youtu.be/jJMJBCGrFw8) (Audiovisual, Spanish, New):35; El desorden que dejas (Subtitled, Spanish, New):41; [Cómo aprender IDIOMAS con la BIBLIA (Sí, EN SERIO)](https://youtu.be/aBOxJ5cRad8) (Audiovisual, Spanish, New):36; "

# Move selection to C18, matching the saved sheet view state
$ws.Range("C18").Select()
